# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (holdings detail, same shape as the
#    other quarterly sheets) right before the "总计" (totals) sheet.
# 2) Insert a new leading row into "总计" summarising the new quarter and
#    shift the existing rows down, renumbering the helper index column.

$wb = $excel.ActiveWorkbook

# --- locate the existing sheets we need as templates / anchors ---------
$totalSheetBefore = $wb.Worksheets.Item(4)   # "总计" (anchor for Add, below)
$templateQtr = $wb.Worksheets.Item(3)        # "2021-Q4" - donor of cell formatting

# --- 1) create the new "2022-Q1" sheet, placed right before "总计" ------
# NOTE: after Worksheets.Add(before), the "before" handle (and any other
# variable bound to that same worksheet) resolves to the newly inserted
# sheet, not the original one any more - so "总计" must be re-fetched by
# position (it is now one slot further along, at index 5) instead of
# reusing $totalSheetBefore for later writes.
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item(5)   # "总计", re-resolved after the Add

# Clone the header (row1, B:H) and body (rows 2-4, A:H) formatting from the
# "2021-Q4" sheet so the new sheet matches the look of its siblings.
$templateQtr.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateQtr.Range("A2:H4").Copy()
$newSheet.Range("A2:H4").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows - columns B:G are stored as text in the source workbook, so
# force a text number-format before assigning them to avoid Excel silently
# coercing the numeric-looking strings into real numbers.
$newSheetData = @(
    @("001481", "华宝油气(QDII)美元",     "39.80", "94.60", "2.31", "0.9194", 4),
    @("162411", "华宝油气(QDII)人民币A",  "39.80", "94.60", "2.31", "0.9194", 4),
    @("007844", "华宝油气(QDII)人民币C",  "12.98", "94.60", "2.31", "0.2998", 4)
)

for ($i = 0; $i -lt $newSheetData.Length; $i++) {
    $row = $i + 2
    $values = $newSheetData[$i]

    $newSheet.Cells.Item($row, 1).Value = $i

    $newSheet.Cells.Item($row, 2).NumberFormat = "@"
    $newSheet.Cells.Item($row, 2).Value = $values[0]
    $newSheet.Cells.Item($row, 2).Style = "Normal"

    $newSheet.Cells.Item($row, 3).NumberFormat = "@"
    $newSheet.Cells.Item($row, 3).Value = $values[1]
    $newSheet.Cells.Item($row, 3).Style = "Normal"

    $newSheet.Cells.Item($row, 4).NumberFormat = "@"
    $newSheet.Cells.Item($row, 4).Value = $values[2]
    $newSheet.Cells.Item($row, 4).Style = "Normal"

    $newSheet.Cells.Item($row, 5).NumberFormat = "@"
    $newSheet.Cells.Item($row, 5).Value = $values[3]
    $newSheet.Cells.Item($row, 5).Style = "Normal"

    $newSheet.Cells.Item($row, 6).NumberFormat = "@"
    $newSheet.Cells.Item($row, 6).Value = $values[4]
    $newSheet.Cells.Item($row, 6).Style = "Normal"

    $newSheet.Cells.Item($row, 7).NumberFormat = "@"
    $newSheet.Cells.Item($row, 7).Value = $values[5]
    $newSheet.Cells.Item($row, 7).Style = "Normal"

    $newSheet.Cells.Item($row, 8).Value = $values[6]
}

# --- 2) insert a new summary row at the top of "总计" -------------------
$totalSheet.Rows.Item(2).Insert()

# The row-insert above copies the header's bold/border formatting down;
# re-stamp the new row with the plain body formatting used by the other
# data rows instead.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 2.14

# Renumber the helper index column (A) for the rows that shifted down
# (rows 3-5 now hold what used to be rows 2-4).
for ($r = 3; $r -le 5; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally active sheet (unrelated to this edit).
$wb.Worksheets.Item(1).Activate()
